# "added 4wk low sales check"
# Updates the forecast figures (MyForecast, Inventory Coverage, Seasonality Index,
# and Reorder Urgency for W23) on the "Forecast Comparison" sheet, and the
# derived totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Row -> [MyForecast(D), InventoryCoverage(H), SeasonalityIndex(L)]
$rows = @(
    @{ Row = 2;  D = 90;  H = 13.03;             L = 0.85 },
    @{ Row = 3;  D = 88;  H = 12.31;             L = 0.98 },
    @{ Row = 4;  D = 90;  H = 11.06;             L = 1.03 },
    @{ Row = 5;  D = 95;  H = 9.529999999999999; L = 1.05 },
    @{ Row = 6;  D = 98;  H = 8.27;              L = 1.17 },
    @{ Row = 7;  D = 100; H = 7.12;              L = 0.9 },
    @{ Row = 8;  D = 101; H = 6.06;              L = 0.8 },
    @{ Row = 9;  D = 100; H = 5.11;              L = 1.12 },
    @{ Row = 10; D = 95;  H = 4.33;              L = 0.8100000000000001 },
    @{ Row = 11; D = 85;  H = 3.72;              L = 0.82 },
    @{ Row = 12; D = 74;  H = 3.12;              L = 0.9399999999999999 },
    @{ Row = 13; D = 65;  H = 2.42;              L = 0.82 },
    @{ Row = 14; D = 61;  H = 1.51;              L = 1.05 },
    @{ Row = 15; D = 59;  H = 0.53;              L = 0.91 },
    @{ Row = 16; D = 53;  H = 0;                 L = 1.02 },
    @{ Row = 17; D = 40;  H = $null;             L = 0.86 }
)

foreach ($r in $rows) {
    $wsForecast.Range("D$($r.Row)").Value = $r.D
    if ($null -ne $r.H) {
        $wsForecast.Range("H$($r.Row)").Value = $r.H
    }
    $wsForecast.Range("L$($r.Row)").Value = $r.L
}

# Row 15 (W23) Reorder Urgency flips from Normal to Urgent due to the
# newly added 4-week low-sales check.
$wsForecast.Range("J15").Value = "Urgent"

# Summary sheet totals, recomputed from the new MyForecast values.
# (Values on this sheet are stored as text, so force text formatting,
# assign, then drop back to the Normal style to avoid leaving a stray
# number format behind.)
$summaryUpdates = @(
    @{ Cell = "B9";  Text = "1294" },
    @{ Cell = "B10"; Text = "762" },
    @{ Cell = "B11"; Text = "363" },
    @{ Cell = "B12"; Text = "101" },
    @{ Cell = "B14"; Text = "40" }
)

foreach ($u in $summaryUpdates) {
    $cell = $wsSummary.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Text
    $cell.Style = "Normal"
}
